# Add files via upload
# Fill in the newly-entered "Actual Result" (column L) values for rows 10-13
# of the MapOfIreland / SearchAPI test cases, and grow row 10 so the wrapped
# text in L10 is fully visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("L10").Value = "dislays covid cases"
$ws.Range("L11").Value = "zooms in and out on map"
$ws.Range("L12").Value = "returns an alert with`nwith covid data"
$ws.Range("L13").Value = "returns empty alert"

# L12's new text wraps onto two lines, so give it a wrapping style distinct
# from its neighbours.
$ws.Range("L12").WrapText = $true

# Row 10 needs to grow to fit the now-taller content.
$ws.Rows.Item(10).RowHeight = 75

# Reflect where the author ended up scrolled to / selecting afterwards.
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("L13").Select()
